$d = $word.ActiveDocument

# Fix accent: "pagina" -> "página"
$d.Content.Find.Execute("de la pagina", $true, $false, $false, $false, $false,
                         $true, 1, $false, "de la página", 2)

# Add two empty paragraphs at the end of the document body
$r = $d.Content
$r.Collapse(0)
$r.InsertXML("<xml/>")

$r2 = $d.Content
$r2.Collapse(0)
$r2.InsertXML("<xml/>")
